$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 38: CSK vs KKR (row 47)
$ws.Range("E47").Value = 60
$ws.Range("H47").Value = 100
$ws.Range("K47").Value = 20
$ws.Range("N47").Value = 40
$ws.Range("Q47").Value = 0
$ws.Range("T47").Value = 80

# Contest 39: RCB vs MI (row 48)
$ws.Range("E48").Value = 100
$ws.Range("H48").Value = 60
$ws.Range("K48").Value = 80
$ws.Range("N48").Value = 40
$ws.Range("Q48").Value = 0
$ws.Range("T48").Value = 20
